# Generate Report for Handoff
# Updates status/date info for the two files that were just handed off
# (c3ac698f-...md and dbb93aef-...md) across the Overview, zh-cn and de-de
# sheets, and widens the "Error Detail" column to fit the new message.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$overviewDate = "2016-09-01 06:29:28"
$zhHandoffDate = "2016-09-01 06:29:24"

$errC3ac698f = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb6911c27ba90083334b46e12d144700b5ade87d/e2e/c3ac698f-71f6-40cd-84e8-3d7ceadc1cbb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d415adcef1289f9c6dfbe36e25f89949ca010b6/e2e/c3ac698f-71f6-40cd-84e8-3d7ceadc1cbb.md."
$errDbb93aef = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb6911c27ba90083334b46e12d144700b5ade87d/e2e/dbb93aef-4adb-4295-86b2-0a7370cfbe80.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d415adcef1289f9c6dfbe36e25f89949ca010b6/e2e/dbb93aef-4adb-4295-86b2-0a7370cfbe80.md."

# ---- Overview sheet: rows for c3ac698f (row 4) and dbb93aef (row 5) ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $readyStatus
$wsOverview.Range("F4").Value = $readyStatus
$wsOverview.Range("G4").Value = $overviewDate

$wsOverview.Range("E5").Value = $readyStatus
$wsOverview.Range("F5").Value = $readyStatus
$wsOverview.Range("G5").Value = $overviewDate

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $readyStatus
$wsZhCn.Range("H4").Value = $zhHandoffDate
$wsZhCn.Range("P4").Value = $errC3ac698f

$wsZhCn.Range("C5").Value = $readyStatus
$wsZhCn.Range("H5").Value = $zhHandoffDate
$wsZhCn.Range("P5").Value = $errDbb93aef

$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $readyStatus
$wsDeDe.Range("H4").Value = $overviewDate
$wsDeDe.Range("P4").Value = $errC3ac698f

$wsDeDe.Range("C5").Value = $readyStatus
$wsDeDe.Range("H5").Value = $overviewDate
$wsDeDe.Range("P5").Value = $errDbb93aef

$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
